$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old transaction rows (3-9); only the header + one "Saldo Inicial" row remain
$ws.Rows("3:9").Delete()

# New header cell E1 ("Saldo Inicial"), reusing the same header formatting as D1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E1").Value = "Saldo Inicial"

# Row 2 now represents the opening balance carried into the next month
$ws.Range("A2").Value = "Saldo Inicial"
$ws.Range("B2").Value = "Inicial"
$ws.Range("C2").Value = -3500
$ws.Range("D2").Value = 45658
$ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = -3500
